# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets: the handoff has
# been handed back and is now in sync with en-US, so each sheet gets a
# "Latest Target File" (F) / "Latest Handback File" (G) pair populated
# (mirroring the existing "Latest Handoff File" / xlf hyperlinks), the
# "Latest Handback DateTime" (H) gets a real timestamp instead of the
# zero-date placeholder, and the shared "Status" text changes from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---- Status text: every cell sharing this string flips to the handback text
$wb.Worksheets.Item("Overview").Range("B2:C3").Replace("Ready for handoff", "Handed back: in sync with en-US")
$wb.Worksheets.Item("zh-cn").Range("C2:C3").Replace("Ready for handoff", "Handed back: in sync with en-US")
$wb.Worksheets.Item("de-de").Range("C2:C3").Replace("Ready for handoff", "Handed back: in sync with en-US")

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (a93456d5-...)
$ws.Range("F2").Value = "a93456d5-0d09-444e-bda1-7b52da4c0df1.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6247fd4845825a836eb16704f32ab8d3a121e9b6/e2e/a93456d5-0d09-444e-bda1-7b52da4c0df1.md", "", "", "a93456d5-0d09-444e-bda1-7b52da4c0df1.md")

$ws.Range("G2").Value = "a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec7dd4af7012b6c289302949df4dff6d73a3eafe/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.zh-cn.xlf", "", "", "a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.zh-cn.xlf")

$ws.Range("H2").Value = "2016-03-23 07:20:18"

# Row 3 (ec6b8f3b-...)
$ws.Range("F3").Value = "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6247fd4845825a836eb16704f32ab8d3a121e9b6/e2e/ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md", "", "", "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md")

$ws.Range("G3").Value = "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec7dd4af7012b6c289302949df4dff6d73a3eafe/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.zh-cn.xlf", "", "", "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.zh-cn.xlf")

$ws.Range("H3").Value = "2016-03-23 07:20:18"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (a93456d5-...)
$ws.Range("F2").Value = "a93456d5-0d09-444e-bda1-7b52da4c0df1.md"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6247fd4845825a836eb16704f32ab8d3a121e9b6/e2e/a93456d5-0d09-444e-bda1-7b52da4c0df1.md", "", "", "a93456d5-0d09-444e-bda1-7b52da4c0df1.md")

$ws.Range("G2").Value = "a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a41cc10490cd2aa17b474087e867df54724454e/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.de-de.xlf", "", "", "a93456d5-0d09-444e-bda1-7b52da4c0df1.8c8a3eb29cadbe10fcab72b5428ae94517749b2d.de-de.xlf")

$ws.Range("H2").Value = "2016-03-23 07:20:34"

# Row 3 (ec6b8f3b-...)
$ws.Range("F3").Value = "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6247fd4845825a836eb16704f32ab8d3a121e9b6/e2e/ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md", "", "", "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.md")

$ws.Range("G3").Value = "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a41cc10490cd2aa17b474087e867df54724454e/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.de-de.xlf", "", "", "ec6b8f3b-b2b3-4dde-b50d-457500a8e69e.6d6da5ff99a296f735830584d35ab6bf71718d45.de-de.xlf")

$ws.Range("H3").Value = "2016-03-23 07:20:34"
